# Weekly update: a new price observation is inserted as the new row 7
# (Madrigal / Primera, week of 2023-06-06), pushing all subsequent rows
# (old rows 7-37) down by one to rows 8-38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 7; this shifts rows 7..37
# down to 8..38 (preserving their data and the date-format style on column D).
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new weekly record.
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C7").Value = 'Arica y Parinacota'
$ws.Range("D7").Value = 45083
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112013
$ws.Range("G7").Value = 'Alcachofa'
$ws.Range("H7").Value = 'Madrigal'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 19000
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = 19500
$ws.Range("N7").Value = '$/caja 40 unidades'
$ws.Range("O7").Value = 'Región de Coquimbo'
$ws.Range("P7").Value = 488
$ws.Range("Q7").Value = 40
$ws.Range("R7").Value = 'Hortaliza'
